$d = $word.ActiveDocument

# Locate the "Website: http://spotabee.buzz/" paragraph (the last paragraph
# that actually holds content) rather than hard-coding its index.
$websitePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*spotabee.buzz*") {
        $websitePara = $candidate
    }
}

# 1. Move the "_GoBack" bookmark from the "...back gardens" paragraph to the
#    end of the "Website: http://spotabee.buzz/" paragraph (after the
#    hyperlink run, before the paragraph mark). Adding a bookmark with a
#    name that already exists on the document relocates it here.
$bookmarkRange = $d.Range($websitePara.Range.End - 2, $websitePara.Range.End - 1)
$bookmarkRange.Bookmarks.Add("_GoBack")

# 2. Remove the trailing empty paragraphs that followed the website
#    paragraph (there were four of them before sectPr).
$tailRange = $d.Range($websitePara.Range.End, $d.Content.End)
$tailRange.Delete()
